$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 to the aggregated DK_Central value
$ws.Range("A2").Value = "c_DK_Central"
$ws.Range("B2").Value = 24959000

# Update row 3 to the aggregated DK_Decentral value
$ws.Range("A3").Value = "c_DK_Decentral"
$ws.Range("B3").Value = 13817000

# Remove the now-obsolete rows 4-7 (DK1/DK2 split rows) entirely
$ws.Range("A4:B7").ClearContents()
